# summer 2022 week 15 games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New matchup rows to append (Player_1, Points_1, Player_2, Points_2)
$newRows = @(
    @(5,2,3,1),
    @(3,1,3,2),
    @(5,0,5,2),
    @(4,1,4,2),
    @(6,1,6,2),
    @(4,2,4,1),
    @(3,0,3,3),
    @(5,2,4,0),
    @(6,2,6,1),
    @(3,1,6,2),
    @(7,2,5,0),
    @(3,1,4,2),
    @(2,0,2,3),
    @(3,2,5,0),
    @(7,3,5,0),
    @(6,0,5,2),
    @(3,0,3,3),
    @(2,2,2,1),
    @(6,3,7,0),
    @(4,0,2,3)
)

$startRow = 1131
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$lastRow = $startRow + $newRows.Count - 1

# Update the view to match where Excel would land after entering this data
$ws.Application.ActiveWindow.ScrollRow = 1127
$ws.Range("H1136").Select()
